$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted into the "Puerro" data block.
# It lands right before the (old) row 114, so every existing row from 114
# through 133 shifts down by one (to 115-134), and the dimension grows from
# A1:R133 to A1:R134.
$ws.Range("A114:R114").EntireRow.Insert()

# Populate the newly inserted row 114 with the new record's data.
$ws.Range("A114").Value = 9
$ws.Range("B114").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C114").Value = "Metropolitana"
$ws.Range("D114").Value = 45154
$ws.Range("E114").Value = 13
$ws.Range("F114").Value = 100112005
$ws.Range("G114").Value = "Puerro"
$ws.Range("H114").Value = "Sin especificar"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 160
$ws.Range("K114").Value = 7000
$ws.Range("L114").Value = 7000
$ws.Range("M114").Value = 7000
$ws.Range("N114").Value = "$/paquete 20 unidades"
$ws.Range("O114").Value = "Provincia de Chacabuco"
$ws.Range("P114").Value = 350
$ws.Range("Q114").Value = 20
$ws.Range("R114").Value = "Hortaliza"
